$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-11, columns B:G (B=ME, C=MAE, D=MSE, E=RMSE, F=SE, G=N)
$data = @(
    @(0.1335670119460798, 0.4099671025390884, 0.3509508125464788, 0.5924110165640734, 0.5841535040152176, 42),
    @(-0.001513401759154079, 0.475841569153532, 0.4177735214472426, 0.6463540217614822, 0.6543817783518403, 41),
    @(0.1309874194782215, 0.4444290441430878, 0.3749385748310419, 0.6123222802014, 0.6057678922250229, 40),
    @(0.02600164070260688, 0.4864429873837449, 0.4097178451191367, 0.6400920598782153, 0.6479243911420136, 39),
    @(0.1413229605261863, 0.4406483824926739, 0.3914688192909578, 0.6256746912661225, 0.6176868865792924, 38),
    @(0.0295561752036937, 0.4627862106754021, 0.4228225359431512, 0.6502480572390441, 0.6585360833922921, 37),
    @(0.1485433660416745, 0.428250065966257, 0.4010703637891622, 0.6333011635779316, 0.6243668417797383, 36),
    @(0.03269697859646408, 0.4719153317718797, 0.4377743913932188, 0.6616452156505168, 0.6704845803472345, 35),
    @(0.08764418001302933, 0.3538399490260787, 0.2806086297826358, 0.5297250511186306, 0.53028072424651, 34),
    @(0.1097826886828203, 0.441331489103254, 0.4022870010623976, 0.6342609881290174, 0.6343733905601003, 33)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
